$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.996.72"
Set-TextValue "E2" "  +5.34%  "
Set-TextValue "D3" "1.879.06"
Set-TextValue "E3" "  +4.00%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "281.99"
Set-TextValue "E6" "  +0.03%  "
Set-TextValue "D7" "0.5253"
Set-TextValue "E7" "  +3.60%  "
Set-TextValue "D8" "0.3537"
Set-TextValue "E8" "  +0.46%  "
Set-TextValue "D9" "45.27"
Set-TextValue "E9" "  +3.66%  "
Set-TextValue "D10" "0.07072"
Set-TextValue "E10" "  +6.69%  "
Set-TextValue "D11" "20.32"
Set-TextValue "E11" "  +1.43%  "
Set-TextValue "D12" "0.8191"
Set-TextValue "E12" "  -2.24%  "
Set-TextValue "D13" "0.07812"
Set-TextValue "E13" "  +0.59%  "
Set-TextValue "D14" "1.889.47"
Set-TextValue "E14" "  +4.58%  "
Set-TextValue "D15" "5.232"
Set-TextValue "E15" "  +2.92%  "
Set-TextValue "E16" "  +3.36%  "
Set-TextValue "E17" "  +0.06%  "
Set-TextValue "D18" "14.59"
Set-TextValue "E18" "  +4.65%  "
Set-TextValue "D19" "0.000008160"
Set-TextValue "E19" "  +2.53%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  +0.07%  "
Set-TextValue "D21" "27.034.09"
Set-TextValue "E21" "  +5.26%  "
Set-TextValue "D22" "4.786"
Set-TextValue "E22" "  +1.26%  "
Set-TextValue "E23" "  +1.77%  "
Set-TextValue "D24" "6.256"
Set-TextValue "E24" "  +3.26%  "
Set-TextValue "D25" "2.404"
Set-TextValue "E25" "  +14.03%  "
Set-TextValue "D26" "147.03"
Set-TextValue "D27" "17.57"
Set-TextValue "E27" "  +3.74%  "
Set-TextValue "D28" "1.667"
Set-TextValue "E28" "  +0.70%  "
Set-TextValue "D29" "113.54"
Set-TextValue "D30" "4.404"
Set-TextValue "E30" "  +1.83%  "
Set-TextValue "D31" "4.394"
Set-TextValue "E31" "  +4.24%  "
Set-TextValue "D32" "0.08889"
Set-TextValue "E32" "  +1.08%  "
Set-TextValue "D33" "0.04916"
Set-TextValue "E33" "  +2.53%  "
Set-TextValue "D34" "1.176"
Set-TextValue "E34" "  +4.38%  "
Set-TextValue "D35" "0.7464"
Set-TextValue "E35" "  +2.89%  "
Set-TextValue "D36" "2.896"
Set-TextValue "E36" "  +1.37%  "
Set-TextValue "D37" "3.291"
Set-TextValue "E37" "  +8.52%  "
Set-TextValue "D38" "2.409"
Set-TextValue "E38" "  +6.06%  "
Set-TextValue "E39" "  +2.61%  "
Set-TextValue "D40" "0.01892"
Set-TextValue "E40" "  +1.60%  "
Set-TextValue "D41" "0.9794"
Set-TextValue "E41" "  +1.76%  "
Set-TextValue "D42" "117.13"
Set-TextValue "E42" "  +2.06%  "
Set-TextValue "D43" "6.327"
Set-TextValue "E43" "  +2.39%  "
Set-TextValue "D44" "8.177"
Set-TextValue "E44" "  +1.81%  "
Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "1.001"
Set-TextValue "E45" "  +0.04%  "
Set-TextValue "B46" "Decentraland"
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.4626"
Set-TextValue "E46" "  +1.25%  "
Set-TextValue "D47" "0.1370"
Set-TextValue "E47" "  -0.90%  "
Set-TextValue "D48" "9.537"
Set-TextValue "E48" "  +3.30%  "
Set-TextValue "D49" "36.85"
Set-TextValue "E49" "  +2.71%  "
Set-TextValue "D50" "1.527"
Set-TextValue "E50" "  +2.17%  "
Set-TextValue "E51" "  +2.36%  "
